$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New party labels in column C
$ws.Range("C6").Value = "A"
$ws.Range("C7").Value = "FRP"
$ws.Range("C8").Value = "H"
$ws.Range("C9").Value = "KRF"
$ws.Range("C10").Value = "MDG"

# Updated vote counts in column D (correlated polling errors)
$ws.Range("D6").Value = 54.2
$ws.Range("D7").Value = 42.5
$ws.Range("D8").Value = 27.4
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 0.8

# Re-apply J formula across the whole range so Excel regroups it as a shared formula
$ws.Range("J6:J12").Formula = "=I6+F6"

# Selection on the sheet
$ws.Range("D11").Select()
